$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025".
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address line "27444 Rondell St, Agoura Hills CA 91301"
#    into two separate lines, "27444 Rondell St" and "Agoura Hills, CA 91301",
#    followed by a new blank line. Embedding carriage returns in the
#    replacement text makes Word clone the paragraph- and run-level
#    formatting of the line being split for each newly created paragraph,
#    which is exactly the formatting (Arial/22, autoSpaceDE/DN off) used in
#    the target document.
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("27444 Rondell St, Agoura Hills CA 91301", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "27444 Rondell St`rAgoura Hills, CA 91301`r", 2) | Out-Null

# 3. Remove the two blank paragraphs that used to sit directly below
#    "...Board of Directors" — a "No Spacing" styled paragraph followed by a
#    "Title" styled paragraph with bold explicitly turned off — while
#    leaving the next "Title" styled paragraph (bold not overridden) intact.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq "Rancho Palma Grande Homeowners Association Board of Directors") {
        $noSpacingPara = $paras.Item($i + 1)
        $titleBoldOffPara = $paras.Item($i + 2)
        # Delete from the bottom up so indices/ranges above stay valid.
        $titleBoldOffPara.Range.Delete()
        $noSpacingPara.Range.Delete()
        break
    }
}
